$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare the two new "Cliente" detail columns (F, G) by copying the
# --- formatting from an existing plain data column, then give the moved
# --- "Usuario Acceso" column (H) the same formatting the old column E had.
$ws.Cells.Item(1,5).Copy() | Out-Null
$ws.Cells.Item(1,6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,7).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,8).PasteSpecial(-4122) | Out-Null

$ws.Range("E2:E5").Copy() | Out-Null
$ws.Range("F2:F5").PasteSpecial(-4122) | Out-Null
$ws.Range("G2:G5").PasteSpecial(-4122) | Out-Null
$ws.Range("H2:H5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Move the old "Nombre Usuario" values (column E) into the new
# --- "Usuario Acceso" column (H) before column E is repurposed.
$ws.Cells.Item(2,8).Value = $ws.Cells.Item(2,5).Value2
$ws.Cells.Item(3,8).Value = "Manuel Espinoza"
$ws.Cells.Item(4,8).Value = $ws.Cells.Item(4,5).Value2
$ws.Cells.Item(5,8).Value = $ws.Cells.Item(5,5).Value2

# --- Header row ---
$ws.Cells.Item(1,1).Value = "ID NIT Acueducto"
$ws.Cells.Item(1,2).Value = "Acueducto"
$ws.Cells.Item(1,3).Value = "Póliza"
$ws.Cells.Item(1,4).Value = "Serial Medidor"
$ws.Cells.Item(1,5).Value = "Cliente"
$ws.Cells.Item(1,6).Value = "Nit/CC Cliente"
$ws.Cells.Item(1,7).Value = "Dirección Cliente"
$ws.Cells.Item(1,8).Value = "Usuario Acceso"

# --- Row 2 (Imaginamos / Supermercado 1) ---
$ws.Cells.Item(2,5).Value = "Supermercado 1"
$ws.Cells.Item(2,6).Value = 800123222
$ws.Cells.Item(2,7).Value = "Calle 3 3 – 56"

# --- Row 3 (Imaginamos / Supermercado 1) ---
$ws.Cells.Item(3,5).Value = "Supermercado 1"
$ws.Cells.Item(3,6).Value = 800123222
$ws.Cells.Item(3,7).Value = "Calle 3 3 – 56"

# --- Row 4 (Globan / Industria 1) ---
$ws.Cells.Item(4,5).Value = "Industria 1"
$ws.Cells.Item(4,6).Value = 900123123
$ws.Cells.Item(4,7).Value = "Calle 30 # 21 26"

# --- Row 5 (Globan / Industria 1) ---
$ws.Cells.Item(5,5).Value = "Industria 1"
$ws.Cells.Item(5,6).Value = 900123123
$ws.Cells.Item(5,7).Value = "Calle 30 # 21 26"

# --- Column widths / row heights to match the new layout ---
$ws.Columns.Item(1).ColumnWidth = 16.593
$ws.Columns.Item(2).ColumnWidth = 10.598
$ws.Columns.Item(3).ColumnWidth = 10.25
$ws.Columns.Item(4).ColumnWidth = 13.756
$ws.Range("E1:H1").ColumnWidth = 15.259

$ws.Rows.Item(1).RowHeight = 13.8
$ws.Rows.Item(4).RowHeight = 13.8

$ws.Range("H4").Select() | Out-Null
